$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.493.56"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").Value = "1.890.52"
$ws.Range("E3").Value = "  +0.86%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.90"
$ws.Range("E5").Value = "  -1.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4705"
$ws.Range("E7").Value = "  -0.50%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2901"
$ws.Range("E8").Value = "  -0.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06500"
$ws.Range("E9").Value = "  +0.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.21"
$ws.Range("E10").Value = "  +0.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07753"
$ws.Range("E11").Value = "  +0.59%  "

$ws.Range("D12").Value = "1.886.34"
$ws.Range("E12").Value = "  +0.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.82"
$ws.Range("E13").Value = "  -0.63%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7265"
$ws.Range("E14").Value = "  -1.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.196"
$ws.Range("E15").Value = "  +0.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "281.68"
$ws.Range("E16").Value = "  +2.92%  "

$ws.Range("D17").Value = "30.487.55"
$ws.Range("E17").Value = "  -0.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.07"

$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007471"
$ws.Range("E20").Value = "  -0.62%  "

$ws.Range("D21").Value = "2.137.49"
$ws.Range("E21").Value = "  +1.04%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.290"
$ws.Range("E22").Value = "  +0.49%  "

$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.268"
$ws.Range("E24").Value = "  +1.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.04"
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.085"
$ws.Range("E26").Value = "  -1.46%  "

$ws.Range("E27").Value = "  +0.90%  "

$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09733"
$ws.Range("E29").Value = "  -2.89%  "

$ws.Range("E30").Value = "  -1.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.470"
$ws.Range("E31").Value = "  -2.49%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.284"
$ws.Range("E32").Value = "  +0.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.150"
$ws.Range("E33").Value = "  +1.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04857"
$ws.Range("E34").Value = "  +1.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.128"
$ws.Range("E35").Value = "  +0.55%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6961"
$ws.Range("E36").Value = "  -0.03%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.717"
$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("E38").Value = "  +2.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.823"
$ws.Range("E39").Value = "  +2.51%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.40"
$ws.Range("E40").Value = "  +3.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.220"
$ws.Range("E41").Value = "  -0.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.994"
$ws.Range("E42").Value = "  +1.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4252"
$ws.Range("E43").Value = "  +1.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8246"
$ws.Range("E45").Value = "  -1.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.48"
$ws.Range("E46").Value = "  -0.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.589"
$ws.Range("E47").Value = "  +2.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.966"
$ws.Range("E48").Value = "  -0.31%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.11"
$ws.Range("E49").Value = "  -0.99%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "913.14"
$ws.Range("E50").Value = "  -0.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05753"
